# A very minimal admin dashboard added
# Applies: Departments typo fix, 2 new Employees, Leaves update + new row,
# Salaries new row, Users role rename + new user, and cursor/active-tab moves.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Departments — fix "Marketting" -> "Marketing"
# ---------------------------------------------------------------------------
$wsDept = $wb.Worksheets.Item("Departments")
$wsDept.Range("B4").Value = "Marketing"

# ---------------------------------------------------------------------------
# 2. Designations — no data change
# ---------------------------------------------------------------------------
$wsDesig = $wb.Worksheets.Item("Designations")

# ---------------------------------------------------------------------------
# 3. Employees — add Ashish Sapkota as employee #2
# ---------------------------------------------------------------------------
$wsEmp = $wb.Worksheets.Item("Employees")
$wsEmp.Range("A3").Value = 2
$wsEmp.Range("B3").Value = "Ashish"
$wsEmp.Range("C3").Value = "Sapkota"
$wsEmp.Range("D3").Value = "ashishspkt6566@gmail.com"
$wsEmp.Range("E3").Value = 19819276566
$wsEmp.Range("F3").Value = "JanashittalPath Tole, Bharatpur-10"
$wsEmp.Range("G3").Value = "Chitwan"
$wsEmp.Range("H3").Value = "Bagmati"
$wsEmp.Range("I3").Value = "Nepal"
$wsEmp.Range("J2").Copy() | Out-Null
$wsEmp.Range("J3").PasteSpecial(-4122) | Out-Null
$wsEmp.Range("J3").Value = 44371
$wsEmp.Range("K3").Value = 4
$wsEmp.Range("L3").Value = 2
$wsEmp.Columns.Item(10).ColumnWidth = 8.72135416666667
$wsEmp.Range("F7").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. Leaves — bump the first leave's day count, add a second leave request
# ---------------------------------------------------------------------------
$wsLeave = $wb.Worksheets.Item("Leaves")
$wsLeave.Range("D2").Value = 5
$wsLeave.Range("A3").Value = 2
$wsLeave.Range("B3").Value = 2
$wsLeave.Range("C2").Copy() | Out-Null
$wsLeave.Range("C3").PasteSpecial(-4122) | Out-Null
$wsLeave.Range("C3").Value = 44357
$wsLeave.Range("D3").Value = 4
$wsLeave.Range("E3").Value = "Due to traffic jam in Chitwan National Park"
$wsLeave.Columns.Item(2).ColumnWidth = 9.94401041666667
$wsLeave.Columns.Item(5).ColumnWidth = 35.6080729166667
$wsLeave.Range("C7").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5. Salaries — add a salary row for employee #2
# ---------------------------------------------------------------------------
$wsSal = $wb.Worksheets.Item("Salaries")
$wsSal.Range("A3").Value = 4
$wsSal.Range("B3").Value = 2
$wsSal.Range("C3").Value = 1
$wsSal.Range("D3").Value = 0
$wsSal.Range("E3").Value = 8500
$wsSal.Range("G11").Select() | Out-Null

# ---------------------------------------------------------------------------
# 6. Vacancies — no data change, becomes the active sheet/tab
# ---------------------------------------------------------------------------
$wsVac = $wb.Worksheets.Item("Vacancies")
$wsVac.Activate() | Out-Null
$wsVac.Range("E20").Select() | Out-Null

# ---------------------------------------------------------------------------
# 7. Users — rename RoleName -> Role, reset Bijay's password, add Ashish's login
# ---------------------------------------------------------------------------
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("C1").Value = "Role"
$wsUsers.Range("D2").Value = "Nopassword"
$wsUsers.Range("A3").Value = 2
$wsUsers.Range("B3").Value = "aacsspkt3@gmail.com"
$wsUsers.Range("C3").Value = 1
$wsUsers.Range("D3").Value = "Password"
$wsUsers.Columns.Item(2).ColumnWidth = 21.9440104166667
$wsUsers.Columns.Item(3).ColumnWidth = 3.83072916666667
$wsUsers.Columns.Item(4).ColumnWidth = 10.3854166666667
$wsUsers.Range("E8").Select() | Out-Null
